# Upload new version with timestamp
# Adds two new pharmacy items ("شامبو هيد اند شولدر400 مل" and
# "شيلز حريمي مزيل عرق") to the item list, right after
# "شامبو نونو 200ملل " (row 84) and before "عضاضه الجو " (old row 85).
# This pushes every following data row down by two rows and bumps the
# running total in the summary row by the sum of the two new prices.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert two blank rows where the new items belong ------------------
$ws.Range("A85:A86").EntireRow.Insert()

# --- 2. Re-use the formatting of the surrounding data rows (style, fonts,
#        borders, …) instead of leaving the freshly inserted rows with the
#        default "no border" style. Row 87 now holds what used to be row
#        85 ("عضاضه الجو "), so its formatting is representative of every
#        other data row in this table.
$ws.Range("A87:N87").Copy()
$ws.Range("A85:N86").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 3. Row heights matching the rest of the table -------------------------
$ws.Rows.Item(85).RowHeight = 25.5
$ws.Rows.Item(86).RowHeight = 24.75

# --- 4. Re-create the merged cells for the two new rows (B:G, H:K, L:M) ---
$ws.Range("B85:G85").Merge()
$ws.Range("H85:K85").Merge()
$ws.Range("L85:M85").Merge()
$ws.Range("B86:G86").Merge()
$ws.Range("H86:K86").Merge()
$ws.Range("L86:M86").Merge()

# --- 5. Fill in the values for the two new rows ----------------------------
$ws.Range("A85").Value = 82
$ws.Range("B85").Value = "شامبو هيد اند شولدر400 مل"
$ws.Range("H85").Value = "0:0"
$ws.Range("L85").Value = 210
$ws.Range("N85").Value = "1:0"

$ws.Range("A86").Value = 83
$ws.Range("B86").Value = "شيلز حريمي مزيل عرق "
$ws.Range("H86").Value = "3:0"
$ws.Range("L86").Value = 75
$ws.Range("N86").Value = "1:0"

# --- 6. Renumber the "م" serial column for every row pushed down ----------
$ws.Range("A87").Value = 84
$ws.Range("A88").Value = 85
$ws.Range("A89").Value = 86
$ws.Range("A90").Value = 87
$ws.Range("A91").Value = 88
$ws.Range("A92").Value = 89
$ws.Range("A93").Value = 90

# --- 7. Fix up row heights so the previously-existing rows keep exactly the
#        same heights they had before the insert (the canonical workbook
#        edits the cell contents of rows 85-91 in place rather than moving
#        whole rows, so every row height below the two brand new rows stays
#        put; only the two new rows 92/93 and the shifted summary/footer
#        rows 94/95 get fresh heights).
$ws.Rows.Item(87).RowHeight = 25.5
$ws.Rows.Item(88).RowHeight = 25.5
$ws.Rows.Item(89).RowHeight = 24.75
$ws.Rows.Item(90).RowHeight = 25.5
$ws.Rows.Item(91).RowHeight = 24.75
$ws.Rows.Item(92).RowHeight = 25.5
$ws.Rows.Item(93).RowHeight = 25.5

# --- 8. Update the running total row (old row 92, now row 94) -------------
$ws.Rows.Item(94).RowHeight = 25.5
$ws.Range("K94").Value = 3717.83

Write-Host "Inserted new pharmacy items and shifted subsequent rows."
